$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.06464670093312194, 0.5485930735930735, 0.6, 0.36, 0.4112690190658179, 0.7, 0.295, 0.4444382947895094, 0.8, 0.26, 0.5776830350197704)
    3 = @(0.1032633878887499, 0.6445833333333333, 0.75, 0.45, 0.4934383013479232, 0.85, 0.4, 0.5337524086809473, 0.85, 0.3325, 0.6548962121472878)
    4 = @(0.01217881055873414, 0.2508333333333334, 0.35, 0.12, 0.1993788955987748, 0.45, 0.095, 0.2504468758412321, 0.45, 0.065, 0.2856531668212503)
    5 = @(0.0007913807527656885, 0.02083333333333333, 0.05, 0.01, 0.01173196815056891, 0.05, 0.005, 0.01173196815056891, 0.1, 0.01, 0.03693028757952595)
    6 = @(0.02648119535087895, 0.3866666666666667, 0.5, 0.19, 0.331422143515167, 0.7, 0.15, 0.379808336442854, 0.8, 0.135, 0.4652837515752271)
    7 = @(0.06306760132418662, 0.4191883116883117, 0.6, 0.24, 0.3019956889389471, 0.8, 0.26, 0.4049221295690922, 0.85, 0.225, 0.5153684564058171)
    8 = @(0.08009889375229304, 0.6273856209150327, 0.75, 0.43, 0.4632436778799587, 0.85, 0.35, 0.4946707556277395, 0.95, 0.29, 0.6545760439140381)
    9 = @(0.1333599610344434, 0.6330555555555556, 0.85, 0.51, 0.5371131665741533, 0.95, 0.48, 0.6088651634217641, 0.95, 0.375, 0.705027786819972)
    10 = @(0.07408224170163608, 0.5176190476190476, 0.7, 0.35, 0.4194558842742613, 0.85, 0.295, 0.4512070000682071, 0.9, 0.26, 0.5902754080035516)
    11 = @(0.1413495868258528, 0.6338095238095238, 0.8, 0.52, 0.5356246564120216, 0.9, 0.46, 0.5968113702335774, 0.9, 0.365, 0.6946581110815202)
    12 = @(0.04780561722605171, 0.3492507002801121, 0.55, 0.21, 0.2551363016939682, 0.7, 0.225, 0.3377482070611812, 0.75, 0.205, 0.4544895806975411)
    13 = @(0.07670534755380193, 0.440021645021645, 0.6, 0.32, 0.3624276790413858, 0.7, 0.28, 0.402292193973452, 0.9, 0.255, 0.5390832808130274)
    14 = @(0.08679691399624502, 0.5867857142857142, 0.8, 0.4400000000000001, 0.5048487178884085, 0.95, 0.39, 0.5896595202628588, 0.95, 0.275, 0.6773282210144329)
    15 = @(0.1366038510539645, 0.5946078431372549, 0.8, 0.45, 0.4622699507445558, 0.95, 0.42, 0.5216825716850432, 1, 0.38, 0.6880541590235631)
    16 = @(0.0007410862354892205, 0.01041666666666667, 0, 0, 0, 0.05, 0.005, 0.01577324383928644, 0.1, 0.005, 0.02928515156065242)
    17 = @(0.001242507213095448, 0.01571895424836601, 0.05, 0.01, 0.01185988563846481, 0.05, 0.01, 0.02278025747937302, 0.15, 0.01, 0.0465413264761115)
    18 = @(0.03033487295889143, 0.3460598633392751, 0.45, 0.15, 0.2347952947152387, 0.6, 0.17, 0.3322877937678121, 0.75, 0.15, 0.4399001810619153)
    19 = @(0.0007997332125239103, 0.01964285714285714, 0.05, 0.01, 0.01320340612862954, 0.1, 0.01, 0.02987007279529621, 0.1, 0.007500000000000001, 0.03842172699868211)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
